$wb = $excel.ActiveWorkbook

# Each of the 29 worksheets (one per backward-elimination step, models "46"
# down to "18") holds a statsmodels OLS summary dumped as one long string in
# cell B2. The script that (re)writes these sheets was re-run after wrapping
# the Excel write in a try/except, which only changed the "Date:"/"Time:"
# stamps baked into that text (everything else is identical). Patch those two
# fields per sheet, in place, without touching anything else in the cell.

$oldTimes = @(
    "16:11:01",
    "16:11:01",
    "16:11:01",
    "16:11:01",
    "16:11:01",
    "16:11:01",
    "16:11:01",
    "16:11:01",
    "16:11:01",
    "16:11:01",
    "16:11:01",
    "16:11:01",
    "16:11:01",
    "16:11:02",
    "16:11:02",
    "16:11:02",
    "16:11:02",
    "16:11:02",
    "16:11:02",
    "16:11:02",
    "16:11:02",
    "16:11:02",
    "16:11:02",
    "16:11:02",
    "16:11:02",
    "16:11:02",
    "16:11:02",
    "16:11:02",
    "16:11:02"
)
$newTimes = @(
    "23:18:40",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41",
    "23:18:41"
)

$oldDate = "Sun, 29 Dec 2019"
$newDate = "Wed, 01 Jan 2020"

for ($i = 0; $i -lt 29; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $cell = $ws.Range("B2")
    $text = [string]$cell.Value2
    $text = $text.Replace($oldDate, $newDate)
    $text = $text.Replace($oldTimes[$i], $newTimes[$i])
    $cell.Value2 = $text
}
